$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 11.5375
$ws.Range("B3").Value = 5.888899999999989
$ws.Range("E5").Value = 12.0961
$ws.Range("B14").Value = 9.109000000000002
$ws.Range("B16").Value = 9.498700000000003
$ws.Range("E16").Value = 13.40300000000001
$ws.Range("B21").Value = 5.787799999999993
$ws.Range("B23").Value = 5.273900000000005
$ws.Range("B25").Value = 5.789599999999991
